$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the LTspice-derived shared-string placeholders in H:I with the actual
# simulation values (Vin/Vout from LTspice), and add a J column comparing them
# to the measured attenuation via 20*LOG10(Vout/Vin).
$hVals = @{}
$iVals = @{}
$hVals[2] = 9.88
$iVals[2] = 9.2
$hVals[3] = 9.86
$iVals[3] = 7.8
$hVals[4] = 9.84
$iVals[4] = 6.29
$hVals[5] = 9.81
$iVals[5] = 4.98
$hVals[6] = 9.8
$iVals[6] = 3.92
$hVals[7] = 9.78
$iVals[7] = 3.03
$hVals[8] = 9.77
$iVals[8] = 2.26
$hVals[9] = 9.76
$iVals[9] = 2.6
$hVals[10] = 9.75
$iVals[10] = 2.05
$hVals[11] = 9.74
$iVals[11] = 0.57
$hVals[12] = 9.74
$iVals[12] = 0.34
$hVals[13] = 9.74
$iVals[13] = 0.26
$hVals[14] = 9.74
$iVals[14] = 0.2
$hVals[15] = 9.74
$iVals[15] = 0.166
$hVals[16] = 9.73
$iVals[16] = 0.124
$hVals[17] = 9.73
$iVals[17] = 0.084
$hVals[18] = 9.73
$iVals[18] = 0.048
$hVals[19] = 9.73
$iVals[19] = 0.038
$hVals[20] = 9.73
$iVals[20] = 0.06
$hVals[21] = 9.73
$iVals[21] = 0.096
$hVals[22] = 9.73
$iVals[22] = 0.292
$hVals[23] = 9.72
$iVals[23] = 0.65
$hVals[24] = 9.71
$iVals[24] = 1
$hVals[25] = 9.67
$iVals[25] = 2.7
$hVals[26] = 9.59
$iVals[26] = 4.59
$hVals[27] = 9.53
$iVals[27] = 5.79
$hVals[28] = 9.48
$iVals[28] = 6.65
$hVals[29] = 9.45
$iVals[29] = 7.2
$hVals[30] = 9.42
$iVals[30] = 7.63
$hVals[31] = 9.4
$iVals[31] = 7.93
$hVals[32] = 9.39
$iVals[32] = 8.16
$hVals[33] = 9.38
$iVals[33] = 8.32

for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 8).Value = $hVals[$r]   # column H
    $ws.Cells.Item($r, 9).Value = $iVals[$r]   # column I
}

$ws.Range("J2").Formula = "=20*LOG10(I2/H2)"
$ws.Range("J3:J33").Formula = "=20*LOG10(I3/H3)"

# Update the saved view/selection state to match the new active cell (K1)
[void]$ws.Range("K1").Select()
